# Add a new "Sheet2" (Budget) worksheet right after the existing "Sheet1",
# make it the active sheet, populate it with the Owner/Budget data and turn
# the range into a table named "Budget" (mirrors the "Vendas" table already
# present on Sheet1).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after Sheet1 so tab order is Sheet1, Sheet2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "Owner"
$ws2.Range("B1").Value = "Budget"

# Data rows
$ws2.Range("A2").Value = "Cosas"
$ws2.Range("B2").Value = 1000
$ws2.Range("A3").Value = "Matheus"
$ws2.Range("B3").Value = 1000
$ws2.Range("A4").Value = "Bruno"
$ws2.Range("B4").Value = 3000

# Turn the populated range into an Excel table named "Budget"
$tbl = $ws2.ListObjects.Add(1, $ws2.Range("A1:B4"), $null, 1)
$tbl.Name = "Budget"

# Make Sheet2 the active/selected sheet (matches activeTab + tabSelected move)
$ws2.Activate()
